$wb = $excel.ActiveWorkbook

# Sheet "ready_to_sale": update Fire row (row 2) in-stock counts
$wsReady = $wb.Worksheets.Item("ready_to_sale")
$wsReady.Range("B2").Value = 28
$wsReady.Range("C2").Value = 72

# Sheet "components": write off materials used (subtract consumed quantities)
$wsComponents = $wb.Worksheets.Item("components")
$wsComponents.Range("B2").Value = 15
$wsComponents.Range("B5").Value = 6
$wsComponents.Range("B6").Value = 4
$wsComponents.Range("B9").Value = 5
$wsComponents.Range("B11").Value = 45
